$wb = $excel.ActiveWorkbook

# Duplicate the ORDERS sheet (same headers/columns/formatting) and place it
# as the new last tab, then rename it to "ETC ORDERS".
$orders = $wb.Worksheets.Item("ORDERS")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$orders.Copy([Type]::Missing, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "ETC ORDERS"

# The copied sheet becomes the active tab (matches activeTab pointing at it).
[void]$newSheet.Select()
[void]$newSheet.Range("G1").Select()
